# Plots.xlsx - plotConfiguration sheet: add an "aggregation" column.
#
# A new column is inserted before the existing "quantiles" column so that
# values placed under it are correctly applied (previously there was no
# dedicated "aggregation" column, so such values ended up misaligned).
#
# Resulting column order on "plotConfiguration":
#   ... I=xAxisLimits, J=yAxisLimits, K=aggregation (new), L=quantiles, M=foldDistance

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("plotConfiguration")

# Insert a new column at K, shifting the old "quantiles"/"foldDistance"
# columns one to the right (K -> L, L -> M).
$null = $ws.Columns("K:K").Insert()

# Header for the newly inserted column.
$ws.Cells.Item(1, 11).Value = "aggregation"

# Match the column's on-disk width to the rest of the table.
$ws.Columns("K:K").ColumnWidth = 10.14

# Leave the selection on the cell that was just edited.
$null = $ws.Range("K2").Select()
